# Update countries & provincias Spain
# Applies the 14-Jun-2020 11:15 data refresh to the "Pais" sheet:
#   - swaps the display names of a few adjacent countries back into
#     alphabetical order (data follows the name it belongs to)
#   - refreshes the day's case/death counters for a handful of countries
#   - bumps the "Datos actualizados..." timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 11:15"

# --- Statistic refresh (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) ---

# Banglades (row 21)
$ws.Range("B21").Value = 87520
$ws.Range("C21").Value = 3141
$ws.Range("D21").Value = 18730
$ws.Range("E21").Value = 67619
$ws.Range("G21").Value = 32
$ws.Range("H21").Value = 1171

# Indonesia (row 34)
$ws.Range("B34").Value = 38277
$ws.Range("C34").Value = 857
$ws.Range("D34").Value = 14531
$ws.Range("E34").Value = 21612
$ws.Range("G34").Value = 43
$ws.Range("H34").Value = 2134

# Polonia (row 40)
$ws.Range("B40").Value = 29392
$ws.Range("C40").Value = 375
$ws.Range("D40").Value = 14226
$ws.Range("E40").Value = 13919
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 1247

# Filipinas (row 41)
$ws.Range("B41").Value = 25930
$ws.Range("C41").Value = 538
$ws.Range("D41").Value = 5954
$ws.Range("E41").Value = 18888
$ws.Range("G41").Value = 14
$ws.Range("H41").Value = 1088

# Barein (row 50)
$ws.Range("E50").Value = 5368
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 41

# Marruecos (row 66)
$ws.Range("B66").Value = 8734
$ws.Range("C66").Value = 42
$ws.Range("D66").Value = 7725
$ws.Range("E66").Value = 797

# El Salvador (row 86)
$ws.Range("B86").Value = 3720
$ws.Range("C86").Value = 117
$ws.Range("D86").Value = 1837
$ws.Range("E86").Value = 1811

# Eslovenia (row 112)
$ws.Range("B112").Value = 1495
$ws.Range("C112").Value = 3
$ws.Range("E112").Value = 27

# Uganda (row 138)
$ws.Range("B138").Value = 696
$ws.Range("C138").Value = 2
$ws.Range("D138").Value = 240
$ws.Range("E138").Value = 456

# --- Reorder countries back into alphabetical order + carry their stats ---

# rows 130/131: Andorra <-> Georgia
$ws.Range("A130").Value = "Georgia"
$ws.Range("B130").Value = 864
$ws.Range("C130").Value = 13
$ws.Range("D130").Value = 703
$ws.Range("E130").Value = 147
$ws.Range("H130").Value = 14

$ws.Range("A131").Value = "Principado de Andorra"
$ws.Range("B131").Value = 853
$ws.Range("D131").Value = 781
$ws.Range("E131").Value = 21
$ws.Range("H131").Value = 51

# rows 206/207: Islas Malvinas <-> Groenlandia (identical stats, name only)
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

# rows 210/211: Seychelles <-> Montserrat
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# rows 213/214: Islas Virgenes Britanicas <-> Papua Nueva Guinea
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
